# Add 2022-Q4 data.
#
# The workbook keeps its newest quarter's fund-holding detail on the sheet
# right after "总计", with each older quarter pushed one tab to the right.
# So applying a new quarter means:
#   - the sheet that used to be "2022-Q3" (with the Q3 detail) is renamed
#     "2022-Q4" and its content is replaced with the new Q4 detail;
#   - a brand-new "2022-Q3" sheet is created after it, holding the Q3 detail
#     that used to live there;
#   - the "总计" summary sheet gets its newest-row-on-top treatment: the old
#     2022-Q3 totals row moves down one row, and the new 2022-Q4 totals are
#     written into the row right under the header.

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)
$wsQ3 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1) Preserve old "2022-Q3" sheet's header/index-column formatting before
#    its content gets replaced, then free up its name.
# ---------------------------------------------------------------------
$wsQ3.Range("B1").Copy()
$wsQ3.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 2) New "2022-Q3" sheet, inserted right after the renamed "2022-Q4" sheet,
#    holding the fund-holding detail that used to be there.
# ---------------------------------------------------------------------
$wsNewQ3 = $wb.Worksheets.Add($null, $wsQ3)
$wsNewQ3.Name = "2022-Q3"

$wsNewQ3.Range("B1:H1").PasteSpecial(-4122)

$wsNewQ3.Range("B1").Value = "基金代码"
$wsNewQ3.Range("C1").Value = "基金名称"
$wsNewQ3.Range("D1").Value = "基金规模"
$wsNewQ3.Range("E1").Value = "股票总仓位"
$wsNewQ3.Range("F1").Value = "仓位占比"
$wsNewQ3.Range("G1").Value = "持有市值(亿元)"
$wsNewQ3.Range("H1").Value = "仓位排名"

$wsNewQ3.Range("A2").Value = 0
$wsNewQ3.Range("B2").Value = "'519772"
$wsNewQ3.Range("C2").Value = "交银新生活力灵活配置混合"
$wsNewQ3.Range("D2").Value = "'49.60"
$wsNewQ3.Range("E2").Value = "'83.55"
$wsNewQ3.Range("F2").Value = "'3.32"
$wsNewQ3.Range("G2").Value = "'1.6467"
$wsNewQ3.Range("H2").Value = 8

$wsNewQ3.Range("A3").Value = 1
$wsNewQ3.Range("B3").Value = "'002137"
$wsNewQ3.Range("C3").Value = "诺安利鑫灵活配置混合A"
$wsNewQ3.Range("D3").Value = "'0.44"
$wsNewQ3.Range("E3").Value = "'76.46"
$wsNewQ3.Range("F3").Value = "'2.42"
$wsNewQ3.Range("G3").Value = "'0.0106"
$wsNewQ3.Range("H3").Value = 10

$wsNewQ3.Range("A4").Value = 2
$wsNewQ3.Range("B4").Value = "'014521"
$wsNewQ3.Range("C4").Value = "诺安利鑫灵活配置混合C"
$wsNewQ3.Range("D4").Value = "'0.05"
$wsNewQ3.Range("E4").Value = "'76.46"
$wsNewQ3.Range("F4").Value = "'2.42"
$wsNewQ3.Range("G4").Value = "'0.0012"
$wsNewQ3.Range("H4").Value = 10

$wsQ3.Range("A2").Copy()
$wsNewQ3.Range("A2:A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3) Replace the renamed "2022-Q4" sheet's content with the new quarter's
#    fund-holding detail. Header/index-column styling matches the
#    bordered/bold/centered style already used on the "总计" sheet.
# ---------------------------------------------------------------------
$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

$wsTotal.Range("B1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)

$wsQ3.Range("A2").Value = 0
$wsQ3.Range("B2").Value = "'002446"
$wsQ3.Range("C2").Value = "广发利鑫灵活配置混合A"
$wsQ3.Range("D2").Value = "'22.53"
$wsQ3.Range("E2").Value = "'73.90"
$wsQ3.Range("F2").Value = "'2.54"
$wsQ3.Range("G2").Value = "'0.5723"
$wsQ3.Range("H2").Value = 7

$wsQ3.Range("A3").Value = 1
$wsQ3.Range("B3").Value = "'011172"
$wsQ3.Range("C3").Value = "广发利鑫灵活配置混合C"
$wsQ3.Range("D3").Value = "'7.03"
$wsQ3.Range("E3").Value = "'73.90"
$wsQ3.Range("F3").Value = "'2.54"
$wsQ3.Range("G3").Value = "'0.1786"
$wsQ3.Range("H3").Value = 7

$wsQ3.Range("A4").Value = 2
$wsQ3.Range("B4").Value = "'002137"
$wsQ3.Range("C4").Value = "诺安利鑫灵活配置混合A"
$wsQ3.Range("D4").Value = "'0.44"
$wsQ3.Range("E4").Value = "'89.87"
$wsQ3.Range("F4").Value = "'3.32"
$wsQ3.Range("G4").Value = "'0.0146"
$wsQ3.Range("H4").Value = 10

$wsQ3.Range("A5").Value = 3
$wsQ3.Range("B5").Value = "'014521"
$wsQ3.Range("C5").Value = "诺安利鑫灵活配置混合C"
$wsQ3.Range("D5").Value = "'0.01"
$wsQ3.Range("E5").Value = "'89.87"
$wsQ3.Range("F5").Value = "'3.32"
$wsQ3.Range("G5").Value = "'0.0003"
$wsQ3.Range("H5").Value = 10

$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2:A5").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4) "总计" sheet: push 2022-Q3 row down, write 2022-Q4 into row 2
# ---------------------------------------------------------------------
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 3
$wsTotal.Range("D3").Value = 1.66

$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 4
$wsTotal.Range("D2").Value = 0.77
